$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.537.53"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").Value = "2.480.37"
$ws.Range("E3").Value = "  -5.89%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'552.92"
$ws.Range("E5").Value = "  -4.81%  "

$ws.Range("D6").Value = "'146.86"
$ws.Range("E6").Value = "  -5.39%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("D9").Value = "2.481.06"
$ws.Range("E9").Value = "  -5.87%  "

$ws.Range("E10").Value = "  -8.40%  "

$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "'5.42"
$ws.Range("E12").Value = "  -6.64%  "

$ws.Range("E13").Value = "  -5.99%  "

$ws.Range("D14").Value = "'26.27"
$ws.Range("E14").Value = "  -7.74%  "

$ws.Range("D15").Value = "2.928.24"
$ws.Range("E15").Value = "  -5.81%  "

$ws.Range("D16").Value = "'0.0000168"
$ws.Range("E16").Value = "  -8.33%  "

$ws.Range("D17").Value = "61.484.25"
$ws.Range("E17").Value = "  -3.54%  "

$ws.Range("D18").Value = "2.482.60"
$ws.Range("E18").Value = "  -5.86%  "

$ws.Range("D19").Value = "'11.17"
$ws.Range("E19").Value = "  -8.00%  "

$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  -8.07%  "

$ws.Range("D21").Value = "'4.22"
$ws.Range("E21").Value = "  -6.81%  "

$ws.Range("D22").Value = "'322.45"
$ws.Range("E22").Value = "  -6.50%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'1.83"
$ws.Range("E24").Value = "  -5.01%  "

$ws.Range("D25").Value = "'64.05"
$ws.Range("E25").Value = "  -5.69%  "

$ws.Range("D26").Value = "0.0₃0990"
$ws.Range("E26").Value = "  -9.13%  "

$ws.Range("D27").Value = "2.609.61"
$ws.Range("E27").Value = "  -5.32%  "

$ws.Range("D28").Value = "'1.53"
$ws.Range("E28").Value = "  -4.94%  "

$ws.Range("D29").Value = "'540.66"
$ws.Range("E29").Value = "  -9.83%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'8.38"
$ws.Range("E31").Value = "  -9.54%  "

$ws.Range("D32").Value = "'7.57"
$ws.Range("E32").Value = "  -6.55%  "

$ws.Range("E33").Value = "  -5.87%  "

$ws.Range("E34").Value = "  -7.64%  "

$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  -8.86%  "

$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = "  -10.53%  "

$ws.Range("E37").Value = "  -10.30%  "

$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "'0.382"
$ws.Range("E39").Value = "  -5.21%  "

$ws.Range("D41").Value = "'149.15"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -8.51%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("E45").Value = "  -8.31%  "

$ws.Range("D46").Value = "'148.06"
$ws.Range("E46").Value = "  -7.44%  "

$ws.Range("E47").Value = "  -7.07%  "

$ws.Range("D48").Value = "'21.02"
$ws.Range("E48").Value = "  -14.78%  "

$ws.Range("E49").Value = "  -8.52%  "

$ws.Range("E50").Value = "  -5.87%  "

$ws.Range("E51").Value = "  -4.99%  "
